# Add "Table IV" and "Table V" worksheets (silhouette-coefficient tables)
# after the existing "Fig 2" sheet, populate their data, and leave
# "Table V" as the active/selected sheet (matching the authored workbook).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table IV
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet, 1, 1)
$ws4.Name = "Table IV"

$ws4.Range("A1").Value = "Cluster"
$ws4.Range("B1").Value = "Silhouette"
$ws4.Range("C1").Value = "Object"

$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = 0.2
$ws4.Range("C2").Value = "a"

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = 0.31
$ws4.Range("C3").Value = "b"

$ws4.Range("A4").Value = 1
$ws4.Range("B4").Value = -0.07
$ws4.Range("C4").Value = "c"

$ws4.Range("A5").Value = 1
$ws4.Range("B5").Value = 0.99
$ws4.Range("C5").Value = "d"

$ws4.Range("A6").Value = 2
$ws4.Range("B6").Value = 0.2
$ws4.Range("C6").Value = "e"

$ws4.Range("A7").Value = 2
$ws4.Range("B7").Value = -0.9
$ws4.Range("C7").Value = "f"

$ws4.Range("A8").Value = 2
$ws4.Range("B8").Value = 0
$ws4.Range("C8").Value = "g"

$ws4.Range("A9").Value = 2
$ws4.Range("B9").Value = 0.23
$ws4.Range("C9").Value = "h"

$ws4.Range("A10").Value = 3
$ws4.Range("B10").Value = 0.1
$ws4.Range("C10").Value = "i"

$ws4.Range("A11").Value = 3
$ws4.Range("B11").Value = 0.4
$ws4.Range("C11").Value = "j"

$ws4.Range("A12").Value = "avg Silhouette"
$ws4.Range("B12").Formula = "=AVERAGE(B2:B11)"

$ws4.Columns.Item(1).ColumnWidth = 16.16
$ws4.Columns.Item(2).ColumnWidth = 12.66

[void]$ws4.Range("A1:C12").Select()

# ---------------------------------------------------------------------
# Table V
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet2, 1, 1)
$ws5.Name = "Table V"

$ws5.Range("A1").Value = "Cluster"
$ws5.Range("B1").Value = "Silhouette"
$ws5.Range("C1").Value = "Object"

$ws5.Range("A2").Value = 1
$ws5.Range("B2").Value = 0.32
$ws5.Range("C2").Value = "a"

$ws5.Range("A3").Value = 1
$ws5.Range("B3").Value = 0.99
$ws5.Range("C3").Value = "b"

$ws5.Range("A4").Value = 1
$ws5.Range("B4").Value = 0.75
$ws5.Range("C4").Value = "c"

$ws5.Range("A5").Value = 1
$ws5.Range("B5").Value = 0.81
$ws5.Range("C5").Value = "d"

$ws5.Range("A6").Value = 1
$ws5.Range("B6").Value = 0.42
$ws5.Range("C6").Value = "e"

$ws5.Range("A7").Value = 2
$ws5.Range("B7").Value = 0.2
$ws5.Range("C7").Value = "f"

$ws5.Range("A8").Value = 2
$ws5.Range("B8").Value = 0.5
$ws5.Range("C8").Value = "g"

$ws5.Range("A9").Value = 2
$ws5.Range("B9").Value = 0.6
$ws5.Range("C9").Value = "h"

$ws5.Range("A10").Value = 2
$ws5.Range("B10").Value = 0.2
$ws5.Range("C10").Value = "i"

$ws5.Range("A11").Value = 2
$ws5.Range("B11").Value = 0.64
$ws5.Range("C11").Value = "j"

$ws5.Range("A12").Value = "avg Silhouette"
$ws5.Range("B12").Formula = "=AVERAGE(B2:B11)"

[void]$ws5.Range("D14").Select()

# "Fig 2" no longer the selected tab -- "Table V" is active/selected now.
$ws5.Activate()

$wb.Save()
